$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 4654.236
$ws.Range("I15").Value = 4654.236
$ws.Range("K15").Value = 13962.708
$ws.Range("M15").Value = -13793.708
$ws.Range("H129").Value = 1097.2727
$ws.Range("I129").Value = 1241.5
$ws.Range("J129").Value = 1043.1875
$ws.Range("K129").Value = 3724.5
$ws.Range("L129").Value = 3129.5625
$ws.Range("M129").Value = 1275.5
$ws.Range("N129").Value = -13129.5625
$ws.Range("H132").Value = 22494.422
$ws.Range("I132").Value = 2932.7026
$ws.Range("J132").Value = 112967.375
$ws.Range("K132").Value = 8798.1078
$ws.Range("L132").Value = 338902.125
$ws.Range("M132").Value = -6268.1078
$ws.Range("N132").Value = -343962.125
$ws.Range("H137").Value = 1484927.9
$ws.Range("I137").Value = 3211498.2
$ws.Range("J137").Value = 5010.4644
$ws.Range("K137").Value = 9634494.600000001
$ws.Range("L137").Value = 15031.3932
$ws.Range("M137").Value = -9631944.600000001
$ws.Range("N137").Value = -20131.3932
$ws.Range("H141").Value = 3699.5557
$ws.Range("I141").Value = 1936.3846
$ws.Range("J141").Value = 8283.8
$ws.Range("K141").Value = 5809.1538
$ws.Range("L141").Value = 24851.4
$ws.Range("M141").Value = -629.1538
$ws.Range("N141").Value = -35211.39999999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 10622.904
$ws.Range("I32").Value = 10067.596
$ws.Range("J32").Value = 13248
$ws.Range("K32").Value = 10067.596
$ws.Range("L32").Value = 13248
$ws.Range("M32").Value = -9780.596
$ws.Range("N32").Value = -13822
$ws.Range("H110").Value = 2077.1
$ws.Range("I110").Value = 2077.1
$ws.Range("J110").Value = 0
$ws.Range("K110").Value = 2077.1
$ws.Range("L110").Value = 0
$ws.Range("M110").Value = -32.09999999999991
$ws.Range("N110").Value = ""
$ws.Range("H132").Value = 13159946
$ws.Range("I132").Value = 16130714
$ws.Range("J132").Value = 3685.1428
$ws.Range("K132").Value = 48392142
$ws.Range("L132").Value = 11055.4284
$ws.Range("M132").Value = -48389612
$ws.Range("N132").Value = -16115.4284

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 3660.077
$ws.Range("I107").Value = 2979.2856
$ws.Range("J107").Value = 4454.3335
$ws.Range("K107").Value = 2979.2856
$ws.Range("L107").Value = 4454.3335
$ws.Range("M107").Value = -1059.2856
$ws.Range("N107").Value = -8294.3335

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1246
$ws.Range("I16").Value = 1287
$ws.Range("K16").Value = 1287
$ws.Range("M16").Value = -1000
$ws.Range("H31").Value = 3626821
$ws.Range("I31").Value = 1927.9166
$ws.Range("J31").Value = 5957109
$ws.Range("K31").Value = 1927.9166
$ws.Range("L31").Value = 5957109
$ws.Range("M31").Value = -1632.9166
$ws.Range("N31").Value = -5957699
$ws.Range("H34").Value = 3626821
$ws.Range("I34").Value = 1927.9166
$ws.Range("J34").Value = 5957109
$ws.Range("K34").Value = 1927.9166
$ws.Range("L34").Value = 5957109
$ws.Range("M34").Value = -1725.9166
$ws.Range("N34").Value = -5957513
$ws.Range("H62").Value = 3504.6365
$ws.Range("I62").Value = 3568.875
$ws.Range("K62").Value = 3568.875
$ws.Range("M62").Value = -2944.875
$ws.Range("H65").Value = 3504.6365
$ws.Range("I65").Value = 3568.875
$ws.Range("K65").Value = 17844.375
$ws.Range("M65").Value = -14724.375
$ws.Range("H113").Value = 1246
$ws.Range("I113").Value = 1287
$ws.Range("K113").Value = 1287
$ws.Range("M113").Value = 883
$ws.Range("H122").Value = 121243.6
$ws.Range("I122").Value = 134559.56
$ws.Range("J122").Value = 1400
$ws.Range("K122").Value = 403678.68
$ws.Range("L122").Value = 4200
$ws.Range("M122").Value = -401228.68
$ws.Range("N122").Value = -9100

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 1485.4459
$ws.Range("J68").Value = 1631.6227
$ws.Range("L68").Value = 4894.8681
$ws.Range("N68").Value = -6516.8681
$ws.Range("H71").Value = 1485.4459
$ws.Range("J71").Value = 1631.6227
$ws.Range("L71").Value = 14684.6043
$ws.Range("N71").Value = -22796.6043
$ws.Range("H86").Value = 1031.25
$ws.Range("I86").Value = 660
$ws.Range("J86").Value = 1650
$ws.Range("K86").Value = 1980
$ws.Range("L86").Value = 4950
$ws.Range("M86").Value = -794
$ws.Range("N86").Value = -7322
$ws.Range("H89").Value = 1031.25
$ws.Range("I89").Value = 660
$ws.Range("J89").Value = 1650
$ws.Range("K89").Value = 5940
$ws.Range("L89").Value = 14850
$ws.Range("M89").Value = -12
$ws.Range("N89").Value = -26706
$ws.Range("H122").Value = 2277.0657
$ws.Range("I122").Value = 487.4898
$ws.Range("J122").Value = 9584.5
$ws.Range("K122").Value = 4387.4082
$ws.Range("L122").Value = 86260.5
$ws.Range("M122").Value = -1937.4082
$ws.Range("N122").Value = -91160.5
$ws.Range("H132").Value = 3419.3076
$ws.Range("I132").Value = 1312
$ws.Range("J132").Value = 4534.9414
$ws.Range("K132").Value = 11808
$ws.Range("L132").Value = 40814.47259999999
$ws.Range("M132").Value = -9278
$ws.Range("N132").Value = -45874.47259999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("I81").Value = 1000
$ws.Range("J81").Value = 10000
$ws.Range("K81").Value = 2000
$ws.Range("L81").Value = 20000
$ws.Range("M81").Value = -939
$ws.Range("N81").Value = -22122
$ws.Range("I84").Value = 1000
$ws.Range("J84").Value = 10000
$ws.Range("K84").Value = 10000
$ws.Range("L84").Value = 100000
$ws.Range("M84").Value = -4696
$ws.Range("N84").Value = -110608
$ws.Range("H122").Value = 2042703.8
$ws.Range("I122").Value = 2382963.5
$ws.Range("J122").Value = 1145.5
$ws.Range("K122").Value = 7148890.5
$ws.Range("L122").Value = 3436.5
$ws.Range("M122").Value = -7146440.5
$ws.Range("N122").Value = -8336.5
$ws.Range("H126").Value = 981711.2
$ws.Range("I126").Value = 1090532.4
$ws.Range("J126").Value = 2320
$ws.Range("K126").Value = 3271597.2
$ws.Range("L126").Value = 6960
$ws.Range("M126").Value = -3269127.2
$ws.Range("N126").Value = -11900
$ws.Range("H132").Value = 989711.25
$ws.Range("I132").Value = 1243241.6
$ws.Range("K132").Value = 3729724.8
$ws.Range("M132").Value = -3727194.8

